$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")

# Rows 2-15 (Z01-Z09, C01-C05): Create/Read/Update/Delete test results all pass -> TRUE
$ws.Range("B2:E15").Value = $true

# Rows 16-24 (D01-D09): Read/Update/Delete test results pass -> TRUE, Create stays FALSE
$ws.Range("C16:E24").Value = $true

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("G31").Select()
